# Update the "Pin Numbers" sheet so column F (PWM Channel component names)
# reflects the corrected wiring as of 02/17/13:
#   - F4 (Shooter Wheel Talon) is no longer wired to that pin; clear it.
#   - F6 is now wired to the Shooter Wheel Talon.
#   - F14/F15/F16 spike assignments are rotated to match the new wiring.
#   - F24/F27 (Passive Hooks L / Fire Launcher) are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# F4 had "Shooter Wheel Talon" - remove it, that pin is now empty.
$ws.Range("F4").ClearContents()

# F6 now carries what used to be in F4.
$ws.Range("F6").Value = "Shooter Wheel Talon"

# F14/F15/F16 rotate: F14<-(old F16), F15<-(old F14), F16<-(old F15)
$ws.Range("F14").Value = "Compressor Spike"
$ws.Range("F15").Value = "Shooter Angle Spike"
$ws.Range("F16").Value = "Arm Spike (Pivot)"

# F24 and F27 swap values.
$ws.Range("F24").Value = "Fire Launcher"
$ws.Range("F27").Value = "Passive Hooks L"

# Update the view: scroll so column D is the leftmost visible column, row 1
# on top, and select F10 as the active cell (matches the saved view state).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F10").Select() | Out-Null
